# Replace deprecated ExtentHTMLReporter for Extent Reporting:
# On the "Sheet2" tab (the active sheet, sheet5.xml), duplicate the
# "TestStepRun" (column I) Yes/No results into a new column J for every
# data row (rows 4-59), mirroring both the value and the cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column I (rows 4-59) into the new column J, carrying over both
# the shared-string value ("Yes"/"No") and the cell formatting/style.
$ws.Range("I4:I59").Copy($ws.Range("J4:J59"))

# Update the view so the newly added column is visible/selected, matching
# the post-edit selection state.
$ws.Range("I4:I59").Select() | Out-Null
